$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 221; this shifts the existing rows
# 221-238 down to 222-239 (data-only shift, as observed in the diff).
$ws.Rows("221:221").Insert()

# Populate the newly inserted row 221 with the new weekly record.
$ws.Range("A221").Value = 11
$ws.Range("B221").Value = "Vega Monumental Concepción"
$ws.Range("C221").Value = "Bíobío"
$ws.Range("D221").Value = 45223
$ws.Range("D221").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E221").Value = 8
$ws.Range("F221").Value = 100112021
$ws.Range("G221").Value = "Ají"
$ws.Range("H221").Value = "Inferno"
$ws.Range("I221").Value = "Primera"
$ws.Range("J221").Value = 15
$ws.Range("K221").Value = 30000
$ws.Range("L221").Value = 30000
$ws.Range("M221").Value = 30000
$ws.Range("N221").Value = "`$/caja 10 kilos"
$ws.Range("O221").Value = "Región de Arica y Parinacota"
$ws.Range("P221").Value = 3000
$ws.Range("Q221").Value = 10
$ws.Range("R221").Value = "Hortaliza"
